$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# "Volume 30   Number  33" -> "...34"
$ws.Range("A8").Value = "Volume 30   Number  34"
# "Report Covering the Week  8/14/2023  Through  8/20/2023" -> new dates
$ws.Range("C9").Value = "Report Covering the Week  8/21/2023  Through  8/27/2023"

# --- Column E width bump --------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 8.699091

# --- Row 14 (Murder) -------------------------------------------------------
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("C14").Value = 2
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 100
$ws.Range("I14").Value = 5
$ws.Range("K14").Value = -58.333333333333
$ws.Range("L14").Value = 25
$ws.Range("M14").Value = -16.666666666666
$ws.Range("N14").Value = -80.769230769230

# --- Row 15 (Rape) ----------------------------------------------------------
$ws.Range("C15").Value = 2
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 23
$ws.Range("J15").Value = 28
$ws.Range("K15").Value = -17.857142857142
$ws.Range("L15").Value = 15
$ws.Range("M15").Value = 21.052631578947
$ws.Range("N15").Value = -42.5

# --- Row 16 (Robbery) -------------------------------------------------------
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 75
$ws.Range("F16").Value = 34
$ws.Range("H16").Value = 30.769230769230
$ws.Range("I16").Value = 300
$ws.Range("J16").Value = 275
$ws.Range("K16").Value = 9.090909090909
$ws.Range("L16").Value = 29.310344827586
$ws.Range("M16").Value = -6.832298136645
$ws.Range("N16").Value = -71.590909090909

# --- Row 17 (Fel. Assault) ---------------------------------------------------
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 87.5
$ws.Range("F17").Value = 49
$ws.Range("G17").Value = 38
$ws.Range("H17").Value = 28.947368421052
$ws.Range("I17").Value = 443
$ws.Range("J17").Value = 405
$ws.Range("K17").Value = 9.382716049382
$ws.Range("L17").Value = 33.033033033033
$ws.Range("M17").Value = 49.662162162162
$ws.Range("N17").Value = 2.546296296296

# --- Row 18 (Burglary) -------------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 150
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 150
$ws.Range("I18").Value = 154
$ws.Range("J18").Value = 134
$ws.Range("K18").Value = 14.925373134328
$ws.Range("L18").Value = 29.411764705882
$ws.Range("M18").Value = -38.888888888888
$ws.Range("N18").Value = -89.821546596166

# --- Row 19 (Gr. Larceny) ----------------------------------------------------
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 122.222222222222
$ws.Range("F19").Value = 71
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 54.347826086956
$ws.Range("I19").Value = 482
$ws.Range("J19").Value = 436
$ws.Range("K19").Value = 10.550458715596
$ws.Range("L19").Value = 11.316397228637
$ws.Range("M19").Value = 35.014005602240
$ws.Range("N19").Value = -17.465753424657

# --- Row 20 (G.L.A.) ---------------------------------------------------------
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 66.666666666666
$ws.Range("F20").Value = 33
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 65
$ws.Range("I20").Value = 227
$ws.Range("J20").Value = 173
$ws.Range("K20").Value = 31.213872832369
$ws.Range("L20").Value = 60.992907801418
$ws.Range("M20").Value = 92.372881355932
$ws.Range("N20").Value = -77.208835341365

# --- Row 21 (TOTAL) ----------------------------------------------------------
$ws.Range("C21").Value = 61
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 103.333333333333
$ws.Range("G21").Value = 141
$ws.Range("H21").Value = 46.808510638297
$ws.Range("I21").Value = 1634
$ws.Range("J21").Value = 1463
$ws.Range("K21").Value = 11.688311688311
$ws.Range("L21").Value = 27.457098283931
$ws.Range("M21").Value = 19.270072992700
$ws.Range("N21").Value = -64.837529588982

# --- Row 22 (Transit) --------------------------------------------------------
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C22").NumberFormat = "General"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = -3.333333333333

# --- Row 24 (Petit Larceny) ---------------------------------------------------
$ws.Range("C24").Value = 48
$ws.Range("D24").Value = 64
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 170
$ws.Range("G24").Value = 236
$ws.Range("H24").Value = -27.966101694915
$ws.Range("I24").Value = 1655
$ws.Range("J24").Value = 1934
$ws.Range("K24").Value = -14.426059979317
$ws.Range("L24").Value = 82.068206820682
$ws.Range("M24").Value = 95.857988165680

# --- Row 25 (Misd. Assault) ---------------------------------------------------
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 77.777777777777
$ws.Range("F25").Value = 62
$ws.Range("G25").Value = 59
$ws.Range("H25").Value = 5.084745762711
$ws.Range("I25").Value = 554
$ws.Range("J25").Value = 520
$ws.Range("K25").Value = 6.538461538461
$ws.Range("L25").Value = 11.244979919678
$ws.Range("M25").Value = -5.299145299145

# --- Row 26 (UCR Rape*) -------------------------------------------------------
$ws.Range("F26").Value = 9
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 40
$ws.Range("J26").Value = 58
$ws.Range("K26").Value = -31.034482758620
$ws.Range("L26").Value = 42.857142857142

# --- Row 27 (Other Sex Crimes) ------------------------------------------------
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 64
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = 28
$ws.Range("L27").Value = -14.666666666666

# --- Row 28 (Shooting Vic.) ----------------------------------------------------
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C28").NumberFormat = "General"
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("L28").Value = -48.148148148148

# --- Row 29 (Shooting Inc.) -----------------------------------------------------
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C29").NumberFormat = "General"
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("L29").Value = -47.619047619047
